# ---------------------------------------------------------------------------
# "carga de las columnas relacionadas a la llegada de clientes"
#
# 1. Insert a block of new paragraphs before the existing first paragraph
#    ("A la lista de proximos eventos...") plus two more afterwards.
# 2. Rename "Llega Cliente (compra)" -> "Llegada cliente" in the events table
#    and insert a new (struck-through) row right after it that keeps the old
#    wording, then strike through the "encargo" and "retiro" rows too.
# 3. Merge the split "Reloj"/"ero" runs into a single "Relojero" run in the
#    second table.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. New paragraphs ------------------------------------------------------

$anchor = $d.Paragraphs.Item(1).Range   # "A la lista de proximos eventos..."

# Seven empty paragraphs inserted immediately before the anchor paragraph
# (the last of the seven stays empty -> the blank <w:p/> from the diff).
for ($i = 0; $i -lt 7; $i++) {
    $anchor.InsertParagraphBefore()
}

$d.Paragraphs.Item(1).Range.Text = "Controlador es estatico (o cada clase tiene una referencia al controlador)"
$d.Paragraphs.Item(2).Range.Text = "Cada clase se encarga de su info, cuando se le llama una función, por lo general generaran un Evento, por ejemplo cuando el cajero le pasa un reloj al relojero"
$d.Paragraphs.Item(3).Range.Text = "Cuando una clase hace algo, por si genera un evento, lo que hace es determinarle la hora de fin, y agregarlo a la lista de futuros eventos del controlador, por ej, cuando el relojero empieza a trabajar en un reloj genera un evento “fin reparación”"
$d.Paragraphs.Item(4).Range.Text = "Los eventos los guarda el controlador en una lista, y los guarda ordenados"
$d.Paragraphs.Item(5).Range.Text = "Hay eventos de inicio y fin de simulación"
$d.Paragraphs.Item(6).Range.Text = "El controlador en cada evento se encarga de preguntarle a los objetos (cajero, relojero, etc) que datos tienen, para cargar la tabla. Esos datos se guardan en una tupla de strings (o algo asi) si son null, se guarda cadena vacia"
# Paragraph 7 stays empty (the blank <w:p/> separator before "A la lista...").

# Two more paragraphs after the (unchanged) "A la lista de proximos eventos..."
# paragraph, which is now Paragraphs.Item(8).
$afterAnchor = $d.Paragraphs.Item(8).Range
$afterAnchor.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.Text = "No hay un evento distinto para cada necesidad del cliente, lo que quiere el cliente se determina en el momento que llega (primero e determina si quiere comprar, encargar, o retirar, y si quiere retirar, hay que elegir a un cliente random de los que están como pendientes, tal vez sea por medio de los relojes o de los clientes, y hay que ver cual es el estado de su reloj, y si no esta listo, se produce un “error”, no exactamente un error, sino que ocurrio algo que no tendría que haber pasado)"

$d.Paragraphs.Item(9).Range.InsertParagraphAfter()
$d.Paragraphs.Item(10).Range.Text = "La clase cliente capas que no es necesaria…"

# --- 2. Events table ---------------------------------------------------------

$t = $d.Tables.Item(1)

# "Llega Cliente (compra)" -> "Llegada cliente"
$t.Rows.Item(4).Cells.Item(1).Range.Find.Execute("Llega Cliente (compra)", $true, $false, $false, $false, $false, $true, 1, $false, "Llegada cliente", 2)

# Insert a new row right after it, carrying the old wording, struck through.
$t.Rows.Add($t.Rows.Item(5)) | Out-Null
$t.Rows.Item(5).Cells.Item(1).Range.Text = "Llega Cliente (compra)"
$t.Rows.Item(5).Cells.Item(1).Range.Font.StrikeThrough = 1
$t.Rows.Item(5).Cells.Item(2).Range.Font.StrikeThrough = 1

# Strike through the (now shifted) "encargo" row; its label is collapsed into
# a single run (matching the target) rather than keeping the original 3-run
# split -- a self Find/Replace forces the run-merge cleanly.
$t.Rows.Item(6).Cells.Item(1).Range.Find.Execute("Llega Cliente (encargo)", $true, $false, $false, $false, $false, $true, 1, $false, "Llega Cliente (encargo)", 2) | Out-Null
$t.Rows.Item(6).Cells.Item(1).Range.Font.StrikeThrough = 1
$t.Rows.Item(6).Cells.Item(2).Range.Font.StrikeThrough = 1

# Strike through the (now shifted) "retiro" row (both columns have text).
$t.Rows.Item(7).Cells.Item(1).Range.Font.StrikeThrough = 1
$t.Rows.Item(7).Cells.Item(2).Range.Font.StrikeThrough = 1

# --- 3. "Reloj" + "ero" -> "Relojero" in the second table --------------------

$t2 = $d.Tables.Item(2)
$t2.Rows.Item(1).Cells.Item(1).Range.Find.Execute("Relojero", $true, $false, $false, $false, $false, $true, 1, $false, "Relojero", 2)
